$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly observations were added to the top of the data block
# (rows 801-871), pushing the existing 71 rows down by two positions
# (old row 801 -> new row 803, ..., old row 871 -> new row 873).
$ws.Range("A801:R802").Insert()

# New row 801: Tomate, Larga vida, Primera - Región de Arica y Parinacota
$ws.Cells.Item(801, 1).Value = 7
$ws.Cells.Item(801, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(801, 3).Value = "Ñuble"
$ws.Cells.Item(801, 4).Value = 45194
$ws.Cells.Item(801, 5).Value = 16
$ws.Cells.Item(801, 6).Value = 100112020
$ws.Cells.Item(801, 7).Value = "Tomate"
$ws.Cells.Item(801, 8).Value = "Larga vida"
$ws.Cells.Item(801, 9).Value = "Primera"
$ws.Cells.Item(801, 10).Value = 300
$ws.Cells.Item(801, 11).Value = 12000
$ws.Cells.Item(801, 12).Value = 12000
$ws.Cells.Item(801, 13).Value = 12000
$ws.Cells.Item(801, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(801, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(801, 16).Value = 667
$ws.Cells.Item(801, 17).Value = 18
$ws.Cells.Item(801, 18).Value = "Hortaliza"

# New row 802: Tomate, Larga vida, Segunda - Región de Arica y Parinacota
$ws.Cells.Item(802, 1).Value = 7
$ws.Cells.Item(802, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(802, 3).Value = "Ñuble"
$ws.Cells.Item(802, 4).Value = 45194
$ws.Cells.Item(802, 5).Value = 16
$ws.Cells.Item(802, 6).Value = 100112020
$ws.Cells.Item(802, 7).Value = "Tomate"
$ws.Cells.Item(802, 8).Value = "Larga vida"
$ws.Cells.Item(802, 9).Value = "Segunda"
$ws.Cells.Item(802, 10).Value = 300
$ws.Cells.Item(802, 11).Value = 10000
$ws.Cells.Item(802, 12).Value = 10000
$ws.Cells.Item(802, 13).Value = 10000
$ws.Cells.Item(802, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(802, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(802, 16).Value = 556
$ws.Cells.Item(802, 17).Value = 18
$ws.Cells.Item(802, 18).Value = "Hortaliza"
